# Update the surface/bottom "Mendota" temperature rows (NSE/RMSE) and
# replace the never-populated "Monona" placeholder rows with the real
# Mendota A4/B1-B4 salt-drop scenario results, adding a new row for B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-5: tidy up the RMSE (column D) rounding ---
$ws.Range("D2").Value = 1.3
$ws.Range("D3").Value = 1.35
$ws.Range("D4").Value = 1.32
$ws.Range("D5").Value = 2.17

# --- Row 6: Mendota / A4 / Constant salt value of 35 ---
$ws.Range("A6").Value = "Mendota"
$ws.Range("B6").Value = "A4"
$ws.Range("C6").Value = 0.81043200000000004
$ws.Range("D6").Value = 2.77
$ws.Range("E6").Value = "Constant salt value of 35"

# --- Fill in the new B1-B4 scenario labels for rows 7-10 ---
$ws.Range("A7").Value = "Mendota"
$ws.Range("B7").Value = "B1"
$ws.Range("A8").Value = "Mendota"
$ws.Range("B8").Value = "B2"
$ws.Range("A9").Value = "Mendota"
$ws.Range("B9").Value = "B3"
$ws.Range("A10").Value = "Mendota"
$ws.Range("B10").Value = "B4"

# --- NSE / RMSE numbers for rows 7-10 ---
$ws.Range("C7").Value = 0.95328029999999997
$ws.Range("D7").Value = 1.37
$ws.Range("C8").Value = 0.95716210000000002
$ws.Range("D8").Value = 1.31
$ws.Range("C9").Value = 0.84670400000000001
$ws.Range("D9").Value = 2.4900000000000002
$ws.Range("C10").Value = 0.45768249999999999
$ws.Range("D10").Value = 4.68

# --- Descriptions, entered row 7, then 9, then 10, then 8 ---
$ws.Range("E7").Value = "Constant salt value of 0.1, then drop to 0 in 2010"
$ws.Range("E9").Value = "Constant salt value of 10, then drop to 0 in 2010"
$ws.Range("E10").Value = "Constant salt value of 35, then drop to 0 in 2010"
$ws.Range("E8").Value = "Constant salt value of 1, then drop to 0 in 2010"

# --- Page setup: explicit portrait orientation (picked up on next print/export) ---
$ws.PageSetup.Orientation = 1
